# "installing the grid properly so that on Safari Browsers it will hopefully
#  be the same as in other Browsers"
#
# Two changes to the "Cross Browser Test" document:
#
#   1. Remove the stray empty paragraph that sits between the Safari-12
#      "Last grid-item won't stretch over whole width of the grid" bullet
#      and the "Firefox 89:" heading.
#
#   2. Add a "Last grid-item won't stretch over whole width of the grid"
#      bullet to the Firefox-89 section too (right after its
#      "Bulletpoints from Unordered List ... list-style: none;" item, and
#      before the trailing blank list paragraph that precedes
#      "Test on Mobile Devices"), mirroring the bullet that already exists
#      for Safari 12 and for iPad Mini 5 / Safari 13.
#
# Paragraphs are located by nearby text rather than hard-coded indices so
# the edit still lands correctly even if unrelated paragraphs shift around.

$d = $word.ActiveDocument

# --- 1. Delete the stray empty paragraph -----------------------------------
# It is the empty paragraph whose previous paragraph is the Safari-12
# "Last grid-item ..." bullet and whose next paragraph is "Firefox 89:".
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $par = $d.Paragraphs($i)
    if ($par.Range.Text.Trim().Length -eq 0) {
        $prevTxt = ""
        if ($i -gt 1) { $prevTxt = $d.Paragraphs($i - 1).Range.Text }
        $nextTxt = ""
        if ($i -lt $d.Paragraphs.Count) { $nextTxt = $d.Paragraphs($i + 1).Range.Text }
        if ($prevTxt -like "*Last grid-item*" -and $nextTxt -like "*Firefox 89*") {
            $par.Range.Delete()
            break
        }
    }
}

# --- 2. Insert the missing bullet under "Firefox 89:" -----------------------
# Find the "Bulletpoints ... list-style: none;" bullet that belongs to the
# "Chrome 91:" heading inside the Firefox-89 section, then add the new
# bullet right after it, inheriting its list formatting (Listenabsatz /
# numbered list / Arial) via InsertParagraphAfter.
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $par = $d.Paragraphs($i)
    if ($par.Range.Text -like "*list-style: none*" -and $i -gt 1 -and `
        $d.Paragraphs($i - 1).Range.Text -like "*Chrome 91*") {

        # Confirm this "Chrome 91:" belongs to the "Firefox 89:" section
        # (walk backwards until we hit either "Firefox 89" or "Safari").
        $inFirefoxSection = $false
        for ($j = $i - 1; $j -ge 1; $j--) {
            $t = $d.Paragraphs($j).Range.Text
            if ($t -like "*Firefox 89*") { $inFirefoxSection = $true; break }
            if ($t -like "*Safari*") { break }
        }

        if ($inFirefoxSection) {
            $par.Range.InsertParagraphAfter()
            $newBullet = $d.Paragraphs($i + 1)
            $newBullet.Range.Text = "Last grid-item won’t stretch over whole width of the grid"
            break
        }
    }
}
